# The annex picture was re-saved as PNG, which nudges its stored size by a
# single EMU (5943600 -> 5943599 EMU while the height stays at 935315 EMU).
# Word exposes shape extents in points (1 pt = 12700 EMU), so drive the
# InlineShape.Width setter with the exact point value that rounds back to
# 5943599 EMU: 5943599 / 12700 = 467.9999212598425 pt.

$d = $word.ActiveDocument

$shp = $d.InlineShapes.Item(1)
$shp.Width = 467.9999212598425
